$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The GSC export rolled forward by one day: the oldest date (2025-10-30) drops off
# and a new day (2026-01-28) is appended, so every remaining rows Date shifts up to
# the next rows previous date. Copy(Destination) preserves the text cell type (so
# Excel does not reinterpret the "yyyy-MM-dd" text as a real date/serial number).
$ws.Range("A3").Copy($ws.Range("A2"))
$ws.Range("A4").Copy($ws.Range("A3"))
$ws.Range("A5").Copy($ws.Range("A4"))
$ws.Range("A6").Copy($ws.Range("A5"))
$ws.Range("A7").Copy($ws.Range("A6"))
$ws.Range("A8").Copy($ws.Range("A7"))
$ws.Range("A9").Copy($ws.Range("A8"))
$ws.Range("A10").Copy($ws.Range("A9"))
$ws.Range("A11").Copy($ws.Range("A10"))
$ws.Range("A12").Copy($ws.Range("A11"))
$ws.Range("A13").Copy($ws.Range("A12"))
$ws.Range("A14").Copy($ws.Range("A13"))
$ws.Range("A15").Copy($ws.Range("A14"))
$ws.Range("A16").Copy($ws.Range("A15"))
$ws.Range("A17").Copy($ws.Range("A16"))
$ws.Range("A18").Copy($ws.Range("A17"))
$ws.Range("A19").Copy($ws.Range("A18"))
$ws.Range("A20").Copy($ws.Range("A19"))
$ws.Range("A21").Copy($ws.Range("A20"))
$ws.Range("A22").Copy($ws.Range("A21"))
$ws.Range("A23").Copy($ws.Range("A22"))
$ws.Range("A24").Copy($ws.Range("A23"))
$ws.Range("A25").Copy($ws.Range("A24"))
$ws.Range("A26").Copy($ws.Range("A25"))
$ws.Range("A27").Copy($ws.Range("A26"))
$ws.Range("A28").Copy($ws.Range("A27"))
$ws.Range("A29").Copy($ws.Range("A28"))
$ws.Range("A30").Copy($ws.Range("A29"))
$ws.Range("A31").Copy($ws.Range("A30"))
$ws.Range("A32").Copy($ws.Range("A31"))
$ws.Range("A33").Copy($ws.Range("A32"))
$ws.Range("A34").Copy($ws.Range("A33"))
$ws.Range("A35").Copy($ws.Range("A34"))
$ws.Range("A36").Copy($ws.Range("A35"))
$ws.Range("A37").Copy($ws.Range("A36"))
$ws.Range("A38").Copy($ws.Range("A37"))
$ws.Range("A39").Copy($ws.Range("A38"))
$ws.Range("A40").Copy($ws.Range("A39"))
$ws.Range("A41").Copy($ws.Range("A40"))
$ws.Range("A42").Copy($ws.Range("A41"))
$ws.Range("A43").Copy($ws.Range("A42"))
$ws.Range("A44").Copy($ws.Range("A43"))
$ws.Range("A45").Copy($ws.Range("A44"))
$ws.Range("A46").Copy($ws.Range("A45"))
$ws.Range("A47").Copy($ws.Range("A46"))
$ws.Range("A48").Copy($ws.Range("A47"))
$ws.Range("A49").Copy($ws.Range("A48"))
$ws.Range("A50").Copy($ws.Range("A49"))
$ws.Range("A51").Copy($ws.Range("A50"))
$ws.Range("A52").Copy($ws.Range("A51"))
$ws.Range("A53").Copy($ws.Range("A52"))
$ws.Range("A54").Copy($ws.Range("A53"))
$ws.Range("A55").Copy($ws.Range("A54"))
$ws.Range("A56").Copy($ws.Range("A55"))
$ws.Range("A57").Copy($ws.Range("A56"))
$ws.Range("A58").Copy($ws.Range("A57"))
$ws.Range("A59").Copy($ws.Range("A58"))
$ws.Range("A60").Copy($ws.Range("A59"))
$ws.Range("A61").Copy($ws.Range("A60"))
$ws.Range("A62").Copy($ws.Range("A61"))
$ws.Range("A63").Copy($ws.Range("A62"))
$ws.Range("A64").Copy($ws.Range("A63"))
$ws.Range("A65").Copy($ws.Range("A64"))
$ws.Range("A66").Copy($ws.Range("A65"))
$ws.Range("A67").Copy($ws.Range("A66"))
$ws.Range("A68").Copy($ws.Range("A67"))
$ws.Range("A69").Copy($ws.Range("A68"))
$ws.Range("A70").Copy($ws.Range("A69"))
$ws.Range("A71").Copy($ws.Range("A70"))
$ws.Range("A72").Copy($ws.Range("A71"))
$ws.Range("A73").Copy($ws.Range("A72"))
$ws.Range("A74").Copy($ws.Range("A73"))
$ws.Range("A75").Copy($ws.Range("A74"))
$ws.Range("A76").Copy($ws.Range("A75"))
$ws.Range("A77").Copy($ws.Range("A76"))
$ws.Range("A78").Copy($ws.Range("A77"))
$ws.Range("A79").Copy($ws.Range("A78"))
$ws.Range("A80").Copy($ws.Range("A79"))
$ws.Range("A81").Copy($ws.Range("A80"))
$ws.Range("A82").Copy($ws.Range("A81"))
$ws.Range("A83").Copy($ws.Range("A82"))
$ws.Range("A84").Copy($ws.Range("A83"))
$ws.Range("A85").Copy($ws.Range("A84"))
$ws.Range("A86").Copy($ws.Range("A85"))
$ws.Range("A87").Copy($ws.Range("A86"))
$ws.Range("A88").Copy($ws.Range("A87"))
$ws.Range("A89").Copy($ws.Range("A88"))
$ws.Range("A90").Copy($ws.Range("A89"))
$ws.Range("A91").Copy($ws.Range("A90"))

# Row 91 becomes the newly appended day. Build it as a text formula result (so it
# stays text, not an auto-converted date) in a scratch cell, then copy it into place.
$ws.Range("Z1").Formula = "=""2026-01-28"""
$ws.Range("Z1").Copy($ws.Range("A91"))
$ws.Range("Z1").Clear()

# Refresh the HTTPS URL counts (column C) with the new export totals.
$ws.Range("C2").Value = 92.0
$ws.Range("C3").Value = 102.0
$ws.Range("C4").Value = 113.0
$ws.Range("C5").Value = 115.0
$ws.Range("C6").Value = 107.0
$ws.Range("C7").Value = 105.0
$ws.Range("C8").Value = 100.0
$ws.Range("C9").Value = 94.0
$ws.Range("C10").Value = 86.0
$ws.Range("C11").Value = 83.0
$ws.Range("C12").Value = 66.0
$ws.Range("C13").Value = 54.0
$ws.Range("C14").Value = 46.0
$ws.Range("C15").Value = 43.0
$ws.Range("C16").Value = 40.0
$ws.Range("C17").Value = 37.0
$ws.Range("C18").Value = 35.0
$ws.Range("C19").Value = 30.0
$ws.Range("C20").Value = 29.0
$ws.Range("C21").Value = 26.0
$ws.Range("C22").Value = 25.0
$ws.Range("C24").Value = 26.0
$ws.Range("C26").Value = 25.0
$ws.Range("C28").Value = 27.0
$ws.Range("C29").Value = 28.0
$ws.Range("C31").Value = 27.0
$ws.Range("C36").Value = 26.0
$ws.Range("C37").Value = 25.0
$ws.Range("C40").Value = 26.0
$ws.Range("C41").Value = 27.0
$ws.Range("C43").Value = 29.0
$ws.Range("C45").Value = 30.0
$ws.Range("C47").Value = 31.0
$ws.Range("C52").Value = 32.0
$ws.Range("C56").Value = 30.0
$ws.Range("C57").Value = 31.0
$ws.Range("C58").Value = 32.0
$ws.Range("C59").Value = 30.0
$ws.Range("C60").Value = 28.0
$ws.Range("C64").Value = 29.0
$ws.Range("C66").Value = 28.0
$ws.Range("C67").Value = 27.0
$ws.Range("C69").Value = 28.0
$ws.Range("C70").Value = 27.0
$ws.Range("C74").Value = 26.0
$ws.Range("C76").Value = 27.0
$ws.Range("C77").Value = 26.0
$ws.Range("C79").Value = 25.0
$ws.Range("C83").Value = 26.0
$ws.Range("C84").Value = 25.0
$ws.Range("C85").Value = 24.0
$ws.Range("C86").Value = 23.0
$ws.Range("C87").Value = 24.0
$ws.Range("C90").Value = 25.0
$ws.Range("C91").Value = 26.0
